# Final output data: append two more rows (12 and 13) of parallel
# memory-usage measurements to the "Memory Usage" sheet, extending the
# existing data table from A1:L11 to A1:L13. Only the cells that have
# actual recorded values are populated (the others are intentionally
# left blank, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("C12").Value = 1285360.0
$ws.Range("D12").Value = 2505376.0
$ws.Range("K12").Value = 7324496.0
$ws.Range("L12").Value = 3682472.0

# Row 13
$ws.Range("C13").Value = 0.0
$ws.Range("D13").Value = 1997600.0

Write-Output "Added rows 12-13 to sheet '$($ws.Name)'"
